# Markov illustrative case ("Power_ThermalGen") is reduced to a single
# generator: the remaining row (the former "VariableGenerator" row, row 8)
# is rescaled and the "BaseLoadGenerator" row (row 9) is removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ScenarioA")

# Rescale MaxProd/Qmax/Qmin for the sole remaining generator (row 8).
$ws.Range("G8").Value = 1000
$ws.Range("I8").Value = 200
$ws.Range("J8").Value = 200

# Drop the BaseLoadGenerator entry (row 9) so only one generator remains.
$ws.Rows.Item(9).Delete()

# Match the saved selection/view state of the edited workbook.
$ws.Range("M8").Select()
